# WordTestCreator.removeValuesFromEachRowRandomCell and tests
#
# The "words_list.xlsx" fixture originally only had 4 data rows (2-5) below
# the header, with the last cell (A5) holding a stray "    " placeholder
# left over from a removed value. This re-fills the sheet down to row 25
# with the repeating Word/Translate pairs (cycling through the three
# maintenance/maintainability/enhance rows) and clears out that stray
# whitespace cell by overwriting it with real data, then moves the
# selection down to reflect where the generator left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three word/translation pairs already present in rows 2-4.
$words = @(
    @("maintenance", "техническое обслуживание"),
    @("maintainability", "ремонтопригодность"),
    @("enhance", "усиливать")
)

# Re-populate rows 5 through 25 (row 5 previously held a leftover blank
# "    " placeholder in column A and nothing in column B) by cycling
# through the three word pairs above.
for ($row = 5; $row -le 25; $row++) {
    $pair = $words[($row - 2) % 3]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}

# Move the view/selection down to where the newly generated rows end.
$ws.Range("A20:B25").Select()

# Scroll the window so row 2 is pinned at the top (topLeftCell moves
# from A1 to A2).
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
